$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column M ("national_id") is inserted after the existing column L ("status").
# Column L (11..12 originally) keeps its 26.5 width; extend the same width group to
# include the new column M (13), matching column L's width exactly.
$ws.Columns.Item(13).ColumnWidth = $ws.Columns.Item(12).ColumnWidth

# Copy the cell formatting (borders/fill/font/number format) from column L into the
# new column M for every existing row (1-10), so the new column visually matches the
# rest of the table (header row style, data row style, and the blank footer rows).
for ($r = 1; $r -le 10; $r++) {
    $srcCell = $ws.Cells.Item($r, 12)
    $dstCell = $ws.Cells.Item($r, 13)
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)
}

# Header cell for the new column.
$ws.Range("M1").Value2 = "national_id"

# Only the second data row has a national_id value in this fixture.
$ws.Range("M2").Value2 = $null
$ws.Range("M3").Value2 = "321"

$excel.CutCopyMode = 0
